# Adding more info in output excel file
#
# This edit touches the "URLs Produtos" sheet:
#  - Column A (row numbering) is filled in completely for rows 3..102
#    (several rows had been left blank before; they now get sequential
#    numbers so every row from 1 to 100 is numbered without gaps).
#  - Two rows (86 and 87) that used a slightly-different "orphan" text
#    style are normalized to use the same style as the rest of column B
#    (Arial 11, black) instead of their own one-off font.
#  - The saved view/selection state is updated to reflect where the
#    author left the cursor (top-left cell back at A1, selection at B114).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("URLs Produtos")

# Make sure this sheet is the active one, same as in the source workbook.
$ws.Activate()

# --- Renumber column A (rows 3-102) so the sequence has no gaps ------------
# Row N (3 <= N <= 102) gets the sequential index N-2 (1, 2, 3, ... 100).
for ($row = 3; $row -le 102; $row++) {
    $ws.Cells.Item($row, 1).Value = $row - 2
}

# --- Normalize the font used on B86/B87 to match the rest of column B -----
$normalizedFontRange = $ws.Range("B86:B87")
$normalizedFontRange.Font.Name = "Arial"
$normalizedFontRange.Font.Size = 11
$normalizedFontRange.Font.Color = 0
$normalizedFontRange.Font.Bold = $false

# --- Update the view/selection state ---------------------------------------
$ws.Range("A1").Select()
$ws.Range("B114").Select()
